# "Generate Report for Handoff"
# Rotate the handoff id (old GUID -> new GUID), the xlf content-hash, and the
# associated handoff timestamps, across the Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

$oldId   = "bf97487f-17b6-4528-b8d9-03a7af2112b9"
$newId   = "6898031e-6338-4b97-b81d-4ac44e82f6ea"
$oldHash = "87bbe253cdc7d08aa6f021fc6f03fb531aa197ac"
$newHash = "7a96ea2c64115eba2d44cd093ab4813d8ba59954"

$newMdName    = "$newId.md"
$newZhXlfName = "$newId.$newHash.zh-cn.xlf"
$newDeXlfName = "$newId.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview": A2 is the handoff-source hyperlink, D2 is the plain
# "latest handoff" timestamp text (no hyperlink).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$ovA2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/406789e9e76a9b4998b0c68b8377d2c220ad627f/e2e/$newMdName"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $ovA2Addr, "", "", $newMdName)

$wsOverview.Range("D2").Value = "2016-12-18 03:12:01"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": A2 (.md hyperlink), B2 (.md extension hyperlink, unchanged),
# D2 (xlf hyperlink), E2 (handoff datetime text, no hyperlink).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhA2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/406789e9e76a9b4998b0c68b8377d2c220ad627f/e2e/$newMdName"
$zhB2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/406789e9e76a9b4998b0c68b8377d2c220ad627f/e2e/$newMdName"
$zhD2Addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b5f568599550851144eec953f8bf6e0e00242c21/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$newZhXlfName"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Addr, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $zhB2Addr, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhD2Addr, "", "", $newZhXlfName)

$wsZh.Range("E2").Value = "2016-03-18 03:11:53"

# ---------------------------------------------------------------------------
# Sheet "de-de": A2 (.md hyperlink), B2 (.md extension hyperlink, unchanged),
# D2 (xlf hyperlink), E2 (handoff datetime text, no hyperlink).
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deA2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/406789e9e76a9b4998b0c68b8377d2c220ad627f/e2e/$newMdName"
$deB2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/406789e9e76a9b4998b0c68b8377d2c220ad627f/e2e/$newMdName"
$deD2Addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f940857ffd32862a1edb5dc1008835bcf0cd349/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$newDeXlfName"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Addr, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $deB2Addr, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deD2Addr, "", "", $newDeXlfName)

$wsDe.Range("E2").Value = "2016-03-18 03:12:01"
